$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header label: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# 2. Remove rows that are no longer present in the final layout.
#    Delete from the bottom up so row numbers of rows not yet processed
#    stay stable.
$ws.Rows(14).Delete()   # "fonte: ibge, ..." footer row
$ws.Rows(8).Delete()    # "grandes regiões" section header row
$ws.Rows(5).Delete()    # "situação do domicílio" section header row

# After the deletions the rows now read (A column labels):
#   1 (title), 2 (column headers), 4 brasil, 5 urbana, 6 rural,
#   7 norte, 8 nordeste, 9 sudeste, 10 sul, 11 centro-oeste

# 3. Write/refresh the numeric data for each remaining data row.
$ws.Range("B5:G5").Value2 = @(0, 5.22, 3.81, 3.05, 3.28, 5.06)
$ws.Range("B6:G6").Value2 = @(0, 9.12, 8.8, 9.92, 13.28, 26.73)
$ws.Range("B7:G7").Value2 = @(0, 10.4, 6.4, 6.95, 9.46, 12.53)
$ws.Range("B8:G8").Value2 = @(0, 6.17, 4.98, 4.5, 6.78, 9.83)
$ws.Range("B9:G9").Value2 = @(0, 10.99, 8.63, 6.06, 5.17, 7.78)
$ws.Range("B10:G10").Value2 = @(0, 14.22, 8.36, 7.01, 6.85, 10.23)
$ws.Range("B11:G11").Value2 = @(0, 18.1, 11.57, 6.04, 7.43, 12.66)

Write-Output "edit applied"
